# Auto-generated edit script applying numeric corrections to Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19, 8).Value2 = 12285.318
$ws.Cells.Item(19, 10).Value2 = 1571.875
$ws.Cells.Item(19, 12).Value2 = 1571.875
$ws.Cells.Item(19, 14).Value2 = -1921.875

# Row 28
$ws.Cells.Item(28, 8).Value2 = 35714824
$ws.Cells.Item(28, 10).Value2 = 956.1
$ws.Cells.Item(28, 12).Value2 = 956.1
$ws.Cells.Item(28, 14).Value2 = -1926.1

# Row 40
$ws.Cells.Item(40, 8).Value2 = 4567.04
$ws.Cells.Item(40, 10).Value2 = 2799
$ws.Cells.Item(40, 12).Value2 = 2799
$ws.Cells.Item(40, 14).Value2 = -3149

# Row 70
$ws.Cells.Item(70, 8).Value2 = 2104.4
$ws.Cells.Item(70, 9).Value2 = 2064
$ws.Cells.Item(70, 11).Value2 = 6192
$ws.Cells.Item(70, 13).Value2 = -5922

# Row 73
$ws.Cells.Item(73, 8).Value2 = 2104.4
$ws.Cells.Item(73, 9).Value2 = 2064
$ws.Cells.Item(73, 11).Value2 = 6192
$ws.Cells.Item(73, 13).Value2 = -5256

# Row 105
$ws.Cells.Item(105, 8).Value2 = 0
$ws.Cells.Item(105, 10).Value2 = 0
$ws.Cells.Item(105, 12).Value2 = 0
$ws.Cells.Item(105, 14).ClearContents()

# Row 107
$ws.Cells.Item(107, 8).Value2 = 17858434
$ws.Cells.Item(107, 9).Value2 = 27778800
$ws.Cells.Item(107, 10).Value2 = 1777.1
$ws.Cells.Item(107, 11).Value2 = 27778800
$ws.Cells.Item(107, 12).Value2 = 1777.1
$ws.Cells.Item(107, 13).Value2 = -27776880
$ws.Cells.Item(107, 14).Value2 = -5617.1

# Row 111
$ws.Cells.Item(111, 8).Value2 = 1465.6
$ws.Cells.Item(111, 9).Value2 = 2164.5
$ws.Cells.Item(111, 10).Value2 = 999.6667
$ws.Cells.Item(111, 11).Value2 = 6493.5
$ws.Cells.Item(111, 12).Value2 = 2999.0001
$ws.Cells.Item(111, 13).Value2 = -3426.5
$ws.Cells.Item(111, 14).Value2 = -9133.000100000001

# Row 112
$ws.Cells.Item(112, 8).Value2 = 993.73914
$ws.Cells.Item(112, 10).Value2 = 993.73914
$ws.Cells.Item(112, 12).Value2 = 2981.21742
$ws.Cells.Item(112, 14).Value2 = -5197.21742

# Row 132
$ws.Cells.Item(132, 8).Value2 = 8946.1
$ws.Cells.Item(132, 9).Value2 = 5782.7896
$ws.Cells.Item(132, 11).Value2 = 17348.3688
$ws.Cells.Item(132, 13).Value2 = -14818.3688

# Row 134
$ws.Cells.Item(134, 8).Value2 = 300000
$ws.Cells.Item(134, 9).Value2 = 0
$ws.Cells.Item(134, 11).Value2 = 0
$ws.Cells.Item(134, 13).ClearContents()

# Row 135
$ws.Cells.Item(135, 8).Value2 = 1806.4286
$ws.Cells.Item(135, 10).Value2 = 3011.6667
$ws.Cells.Item(135, 12).Value2 = 27105.0003
$ws.Cells.Item(135, 14).Value2 = -32175.0003

# Row 137
$ws.Cells.Item(137, 8).Value2 = 1005.72095
$ws.Cells.Item(137, 9).Value2 = 1112.56
$ws.Cells.Item(137, 10).Value2 = 857.3333
$ws.Cells.Item(137, 11).Value2 = 3337.68
$ws.Cells.Item(137, 12).Value2 = 2571.9999
$ws.Cells.Item(137, 13).Value2 = -787.6799999999998
$ws.Cells.Item(137, 14).Value2 = -7671.9999

# Row 138
$ws.Cells.Item(138, 8).Value2 = 1447.75
$ws.Cells.Item(138, 9).Value2 = 1447.75
$ws.Cells.Item(138, 10).Value2 = 0
$ws.Cells.Item(138, 11).Value2 = 4343.25
$ws.Cells.Item(138, 12).Value2 = 0
$ws.Cells.Item(138, 13).Value2 = 796.75
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 133
$ws.Cells.Item(133, 8).Value2 = 50000
$ws.Cells.Item(133, 10).Value2 = 50000
$ws.Cells.Item(133, 12).Value2 = 50000
$ws.Cells.Item(133, 14).Value2 = -55060

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Cells.Item(5, 8).Value2 = 2000
$ws.Cells.Item(5, 9).Value2 = 0
$ws.Cells.Item(5, 10).Value2 = 2000
$ws.Cells.Item(5, 11).Value2 = 0
$ws.Cells.Item(5, 12).Value2 = 2000
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).Value2 = -2226

# Row 14
$ws.Cells.Item(14, 8).Value2 = 1944.5
$ws.Cells.Item(14, 10).Value2 = 1999
$ws.Cells.Item(14, 12).Value2 = 1999
$ws.Cells.Item(14, 14).Value2 = -2343

# Row 33
$ws.Cells.Item(33, 8).Value2 = 621
$ws.Cells.Item(33, 9).Value2 = 621
$ws.Cells.Item(33, 11).Value2 = 621
$ws.Cells.Item(33, 13).Value2 = -285

# Row 80
$ws.Cells.Item(80, 8).Value2 = 840.1539
$ws.Cells.Item(80, 9).Value2 = 654.5
$ws.Cells.Item(80, 10).Value2 = 922.6667
$ws.Cells.Item(80, 11).Value2 = 654.5
$ws.Cells.Item(80, 12).Value2 = 922.6667
$ws.Cells.Item(80, 13).Value2 = 343.5
$ws.Cells.Item(80, 14).Value2 = -2918.6667

# Row 83
$ws.Cells.Item(83, 8).Value2 = 840.1539
$ws.Cells.Item(83, 9).Value2 = 654.5
$ws.Cells.Item(83, 10).Value2 = 922.6667
$ws.Cells.Item(83, 11).Value2 = 3272.5
$ws.Cells.Item(83, 12).Value2 = 4613.3335
$ws.Cells.Item(83, 13).Value2 = 1719.5
$ws.Cells.Item(83, 14).Value2 = -14597.3335

# Row 107
$ws.Cells.Item(107, 8).Value2 = 4303.1377
$ws.Cells.Item(107, 9).Value2 = 4056.074
$ws.Cells.Item(107, 10).Value2 = 4518.3228
$ws.Cells.Item(107, 11).Value2 = 4056.074
$ws.Cells.Item(107, 12).Value2 = 4518.3228
$ws.Cells.Item(107, 13).Value2 = -2136.074
$ws.Cells.Item(107, 14).Value2 = -8358.3228

# Row 134
$ws.Cells.Item(134, 8).Value2 = 5221.241
$ws.Cells.Item(134, 9).Value2 = 4939.115
$ws.Cells.Item(134, 11).Value2 = 14817.345
$ws.Cells.Item(134, 13).Value2 = -12282.345

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value2 = 2398.6
$ws.Cells.Item(16, 9).Value2 = 1973
$ws.Cells.Item(16, 11).Value2 = 1973
$ws.Cells.Item(16, 13).Value2 = -1686

# Row 32
$ws.Cells.Item(32, 8).Value2 = 900
$ws.Cells.Item(32, 9).Value2 = 900
$ws.Cells.Item(32, 11).Value2 = 900
$ws.Cells.Item(32, 13).Value2 = -584

# Row 59
$ws.Cells.Item(59, 8).Value2 = 37250
$ws.Cells.Item(59, 9).Value2 = 49500
$ws.Cells.Item(59, 10).Value2 = 25000
$ws.Cells.Item(59, 11).Value2 = 49500
$ws.Cells.Item(59, 12).Value2 = 25000
$ws.Cells.Item(59, 13).Value2 = -48355
$ws.Cells.Item(59, 14).Value2 = -27290

# Row 64
$ws.Cells.Item(64, 8).Value2 = 33999.5
$ws.Cells.Item(64, 9).Value2 = 8999
$ws.Cells.Item(64, 10).Value2 = 59000
$ws.Cells.Item(64, 11).Value2 = 8999
$ws.Cells.Item(64, 12).Value2 = 59000
$ws.Cells.Item(64, 13).Value2 = -8751
$ws.Cells.Item(64, 14).Value2 = -59496

# Row 67
$ws.Cells.Item(67, 8).Value2 = 33999.5
$ws.Cells.Item(67, 9).Value2 = 8999
$ws.Cells.Item(67, 10).Value2 = 59000
$ws.Cells.Item(67, 11).Value2 = 8999
$ws.Cells.Item(67, 12).Value2 = 59000
$ws.Cells.Item(67, 13).Value2 = -8141
$ws.Cells.Item(67, 14).Value2 = -60716

# Row 103
$ws.Cells.Item(103, 8).Value2 = 20864.666
$ws.Cells.Item(103, 9).Value2 = 20864.666
$ws.Cells.Item(103, 11).Value2 = 20864.666
$ws.Cells.Item(103, 13).Value2 = -19692.666

# Row 113
$ws.Cells.Item(113, 8).Value2 = 2398.6
$ws.Cells.Item(113, 9).Value2 = 1973
$ws.Cells.Item(113, 11).Value2 = 1973
$ws.Cells.Item(113, 13).Value2 = 197

# Row 134
$ws.Cells.Item(134, 8).Value2 = 8232.75
$ws.Cells.Item(134, 9).Value2 = 9243.125
$ws.Cells.Item(134, 10).Value2 = 5201.625
$ws.Cells.Item(134, 11).Value2 = 27729.375
$ws.Cells.Item(134, 12).Value2 = 15604.875
$ws.Cells.Item(134, 13).Value2 = -25194.375
$ws.Cells.Item(134, 14).Value2 = -20674.875

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Cells.Item(14, 8).Value2 = 248.11111
$ws.Cells.Item(14, 9).Value2 = 248.11111
$ws.Cells.Item(14, 11).Value2 = 744.3333299999999
$ws.Cells.Item(14, 13).Value2 = -571.3333299999999

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Cells.Item(113, 8).Value2 = 3110.3076
$ws.Cells.Item(113, 9).Value2 = 2744.6667
$ws.Cells.Item(113, 10).Value2 = 3933
$ws.Cells.Item(113, 11).Value2 = 2744.6667
$ws.Cells.Item(113, 12).Value2 = 3933
$ws.Cells.Item(113, 13).Value2 = -574.6667000000002
$ws.Cells.Item(113, 14).Value2 = -8273

# Row 132
$ws.Cells.Item(132, 8).Value2 = 4021.2703
$ws.Cells.Item(132, 9).Value2 = 3899.6562
$ws.Cells.Item(132, 10).Value2 = 4799.6
$ws.Cells.Item(132, 11).Value2 = 11698.9686
$ws.Cells.Item(132, 12).Value2 = 14398.8
$ws.Cells.Item(132, 13).Value2 = -9168.9686
$ws.Cells.Item(132, 14).Value2 = -19458.8

$ws = $wb.Worksheets.Item("LTW")
# Row 10
$ws.Cells.Item(10, 8).Value2 = 599.4
$ws.Cells.Item(10, 9).Value2 = 496.5
$ws.Cells.Item(10, 11).Value2 = 496.5
$ws.Cells.Item(10, 13).Value2 = -356.5

# Row 61
$ws.Cells.Item(61, 8).Value2 = 14057
$ws.Cells.Item(61, 9).Value2 = 25988.25
$ws.Cells.Item(61, 10).Value2 = 2125.75
$ws.Cells.Item(61, 11).Value2 = 25988.25
$ws.Cells.Item(61, 12).Value2 = 2125.75
$ws.Cells.Item(61, 13).Value2 = -25786.25
$ws.Cells.Item(61, 14).Value2 = -2529.75

# Row 106
$ws.Cells.Item(106, 8).Value2 = 14455.833
$ws.Cells.Item(106, 10).Value2 = 14455.833
$ws.Cells.Item(106, 12).Value2 = 14455.833
$ws.Cells.Item(106, 14).Value2 = -16979.833

# Row 113
$ws.Cells.Item(113, 8).Value2 = 14057
$ws.Cells.Item(113, 9).Value2 = 25988.25
$ws.Cells.Item(113, 10).Value2 = 2125.75
$ws.Cells.Item(113, 11).Value2 = 25988.25
$ws.Cells.Item(113, 12).Value2 = 2125.75
$ws.Cells.Item(113, 13).Value2 = -23818.25
$ws.Cells.Item(113, 14).Value2 = -6465.75

# Row 132
$ws.Cells.Item(132, 8).Value2 = 61978.6
$ws.Cells.Item(132, 9).Value2 = 61978.6
$ws.Cells.Item(132, 11).Value2 = 185935.8
$ws.Cells.Item(132, 13).Value2 = -183405.8

$ws = $wb.Worksheets.Item("WVR")
# Row 38
$ws.Cells.Item(38, 8).Value2 = 24963.334
$ws.Cells.Item(38, 9).Value2 = 24900
$ws.Cells.Item(38, 10).Value2 = 24995
$ws.Cells.Item(38, 11).Value2 = 24900
$ws.Cells.Item(38, 12).Value2 = 24995
$ws.Cells.Item(38, 13).Value2 = -24427
$ws.Cells.Item(38, 14).Value2 = -25941

# Row 56
$ws.Cells.Item(56, 8).Value2 = 23533
$ws.Cells.Item(56, 9).Value2 = 10285
$ws.Cells.Item(56, 10).Value2 = 30157
$ws.Cells.Item(56, 11).Value2 = 10285
$ws.Cells.Item(56, 12).Value2 = 30157
$ws.Cells.Item(56, 13).Value2 = -9571
$ws.Cells.Item(56, 14).Value2 = -31585

# Row 107
$ws.Cells.Item(107, 8).Value2 = 1026.2
$ws.Cells.Item(107, 9).Value2 = 1438.2222
$ws.Cells.Item(107, 10).Value2 = 408.16666
$ws.Cells.Item(107, 11).Value2 = 4314.6666
$ws.Cells.Item(107, 12).Value2 = 1224.49998
$ws.Cells.Item(107, 13).Value2 = -2394.6666
$ws.Cells.Item(107, 14).Value2 = -5064.499980000001

# Row 132
$ws.Cells.Item(132, 8).Value2 = 3441.0667
$ws.Cells.Item(132, 9).Value2 = 3519.1072
$ws.Cells.Item(132, 11).Value2 = 10557.3216
$ws.Cells.Item(132, 13).Value2 = -8027.321599999999

# Row 133
$ws.Cells.Item(133, 8).Value2 = 37290.332
$ws.Cells.Item(133, 10).Value2 = 37290.332
$ws.Cells.Item(133, 12).Value2 = 37290.332
$ws.Cells.Item(133, 14).Value2 = -47410.332

# Row 136
$ws.Cells.Item(136, 8).Value2 = 2756.5557
$ws.Cells.Item(136, 9).Value2 = 2830.2856
$ws.Cells.Item(136, 11).Value2 = 8490.856800000001
$ws.Cells.Item(136, 13).Value2 = -5940.856800000001

